$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.445.02"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.558.69"
$ws.Range("E3").Value = "  +2.04%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "619.91"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.69"
$ws.Range("E6").Value = "  +5.03%  "
$ws.Range("D7").Value = "3.557.05"
$ws.Range("E7").Value = "  +2.04%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.42%  "
$ws.Range("E10").Value = "  +6.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.35"
$ws.Range("E11").Value = "  +6.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.440"
$ws.Range("E12").Value = "  +4.34%  "
$ws.Range("E13").Value = "  +3.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.29"
$ws.Range("E14").Value = "  +6.07%  "
$ws.Range("D15").Value = "4.162.45"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.560.86"
$ws.Range("E16").Value = "  +2.18%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "68.688.60"
$ws.Range("E17").Value = "  +2.35%  "
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  +5.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.99"
$ws.Range("E20").Value = "  +6.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.07"
$ws.Range("E21").Value = "  +11.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "458.81"
$ws.Range("E22").Value = "  +2.69%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.643"
$ws.Range("E23").Value = "  +3.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.62"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.55"
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("D27").Value = "3.704.96"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.07"
$ws.Range("E29").Value = "  +9.49%  "
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("E31").Value = "  +8.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.171"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("E34").Value = "  +5.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.26"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("E36").Value = "  +4.77%  "
$ws.Range("D37").Value = "3.553.44"
$ws.Range("E37").Value = "  +2.10%  "
$ws.Range("E38").Value = "  +4.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.38"
$ws.Range("E39").Value = "  +8.94%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.62"
$ws.Range("E41").Value = "  +4.67%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0920"
$ws.Range("E42").Value = "  +5.62%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E44").Value = "  +3.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "30.56"
$ws.Range("E45").Value = "  +15.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.902"
$ws.Range("E46").Value = "  +2.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.35"
$ws.Range("E47").Value = "  +7.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.79"
$ws.Range("E48").Value = "  +0.94%  "
$ws.Range("E49").Value = "  +6.08%  "
$ws.Range("E50").Value = "  +3.67%  "
$ws.Range("E51").Value = "  +8.55%  "
